# iApply_Integration.xlsx — "adding RDW file intergration to repo"
#
# 1. sheet "sp_view_prod_updation" (index 2): add the RDW numbers block
#    (E45 helper formula, E49:F53 labelled totals, F54/F55 rollup formulas).
# 2. sheet "iApply_Integration_Overall_Stat" (index 1): row 13 gets real
#    status values instead of "-", the whole status-filter is cleared so
#    every previously-filtered-out row becomes visible again, and a second
#    (smaller) RDW numbers block is appended at rows 56-60.
#
# NOTE on ordering: the new one-letter labels ("r","m","f","t","g") must be
# entered on sheet 2 *before* sheet 1 so they land in the shared-string table
# in that exact order (matches the workbook's own save order).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) sp_view_prod_updation (sheet 2) — new rows 45, 49-55
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("E45").Formula = "=10.4*4"

$ws2.Range("E49").Value = "r"
$ws2.Range("F49").Value = 660
$ws2.Range("E50").Value = "m"
$ws2.Range("F50").Value = 35
$ws2.Range("E51").Value = "f"
$ws2.Range("F51").Value = 500
$ws2.Range("E52").Value = "t"
$ws2.Range("F52").Value = 42
$ws2.Range("E53").Value = "g"
$ws2.Range("F53").Value = 73

$ws2.Range("F54").Formula = "=F53+F52+F51+F50+F49"
$ws2.Range("F55").Formula = "=4400-1310"

# ---------------------------------------------------------------------
# 2) iApply_Integration_Overall_Stat (sheet 1)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# -- row 13 used to be all "-" placeholders; fill in the real statuses
$ws1.Range("D13").Value = "COMPLETE"
$ws1.Range("E13").Value = "COMPLETE"
$ws1.Range("F13").Value = "COMPLETE"
$ws1.Range("G13").Value = "NOT REQUIRED"
$ws1.Range("H13").Value = "NOT REQUIRED"
$ws1.Range("I13").Value = "LSA ONLY [ COMPLETED]"
$ws1.Range("J13").Value = "COMPLETE"
$ws1.Range("K13").Value = "NOT REQUIRED"

# row 13 is taller than default (merged-look wrapped text); restore that
# explicitly since writing the cells above can reset the autofit height.
$ws1.Rows.Item(13).RowHeight = 25.5

# -- new RDW numbers block at the bottom of the sheet
$ws1.Range("D56").Value = "r"
$ws1.Range("E56").Value = 660
$ws1.Range("D57").Value = "f"
$ws1.Range("E57").Value = 500
$ws1.Range("D58").Value = "m"
$ws1.Range("E58").Value = 35
$ws1.Range("D59").Value = "t"
$ws1.Range("E59").Formula = "=28*4"
$ws1.Range("E60").Formula = "=E59+E58+E57+E56"

# -- the sheet was filtered down to COMPLETE/NOT REQUIRED rows only;
#    clear that filter so every row is visible again
$ws1.ShowAllData()
$ws1.Rows("1:49").Hidden = $false

# ---------------------------------------------------------------------
# 3) selections, matching the saved cursor position in each sheet
# ---------------------------------------------------------------------
$ws2.Range("F55").Select()
$ws1.Range("A13").Select()
